$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rich-text partial edits (Volume/Number and Report-Covering-the-Week cells) ---
# Edit the later substring first so earlier character offsets remain valid.
$ws.Range("A8").Characters(21,2).Text = "49"
$ws.Range("C9").Characters(48,9).Text = "12/11/2022"
$ws.Range("C9").Characters(27,10).Text = "12/5/2022"

# --- Cells switching from a plain number to the text placeholder "0"/"***.*" (style 14) ---
function Set-TextPlaceholder($dstAddr, $srcAddr) {
    $ws.Range($dstAddr).Value = "__tmp__"
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4104) | Out-Null
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}
Set-TextPlaceholder "C17" "C14"
Set-TextPlaceholder "D18" "C14"
Set-TextPlaceholder "E18" "E14"

# --- Cells switching from the text placeholder to a plain number (style 15) ---
function Set-NumberFromPlaceholder($dstAddr, $srcStyleAddr, $value) {
    $ws.Range($dstAddr).Value = $value
    $ws.Range($srcStyleAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}
Set-NumberFromPlaceholder "C23" "I14" 1
Set-NumberFromPlaceholder "C26" "I14" 1

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("L15").Value = 14.285714285714
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 62
$ws.Range("J16").Value = 42
$ws.Range("K16").Value = 47.619047619047
$ws.Range("L16").Value = 29.166666666666
$ws.Range("M16").Value = -43.636363636363
$ws.Range("N16").Value = -79.124579124579
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -30
$ws.Range("I17").Value = 120
$ws.Range("J17").Value = 118
$ws.Range("K17").Value = 1.694915254237
$ws.Range("L17").Value = 15.384615384615
$ws.Range("M17").Value = -7.692307692307
$ws.Range("N17").Value = -58.762886597938
$ws.Range("C18").Value = 1
$ws.Range("I18").Value = 77
$ws.Range("K18").Value = 14.925373134328
$ws.Range("L18").Value = -14.444444444444
$ws.Range("M18").Value = -61.306532663316
$ws.Range("N18").Value = -94.017094017094
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -20.512820512820
$ws.Range("I19").Value = 346
$ws.Range("J19").Value = 319
$ws.Range("K19").Value = 8.463949843260
$ws.Range("L19").Value = 49.137931034482
$ws.Range("M19").Value = -11.053984575835
$ws.Range("N19").Value = -56.803995006242
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 300
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 130
$ws.Range("J20").Value = 76
$ws.Range("K20").Value = 71.052631578947
$ws.Range("L20").Value = 94.029850746268
$ws.Range("M20").Value = 16.071428571428
$ws.Range("N20").Value = -95.107263831388
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 16.666666666666
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = -4.615384615384
$ws.Range("I21").Value = 744
$ws.Range("J21").Value = 631
$ws.Range("K21").Value = 17.908082408874
$ws.Range("L21").Value = 35.519125683060
$ws.Range("M21").Value = -22.5
$ws.Range("N21").Value = -86.103847590586
$ws.Range("F23").Value = 3
$ws.Range("I23").Value = 20
$ws.Range("K23").Value = -42.857142857142
$ws.Range("L23").Value = 11.111111111111
$ws.Range("M23").Value = -9.090909090909
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 118.181818181818
$ws.Range("F24").Value = 71
$ws.Range("G24").Value = 47
$ws.Range("H24").Value = 51.063829787234
$ws.Range("I24").Value = 824
$ws.Range("J24").Value = 479
$ws.Range("K24").Value = 72.025052192066
$ws.Range("L24").Value = 36.423841059602
$ws.Range("M24").Value = -49.167180752621
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -83.333333333333
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 18.181818181818
$ws.Range("J25").Value = 240
$ws.Range("K25").Value = 32.5
$ws.Range("L25").Value = 43.891402714932
$ws.Range("M25").Value = -37.890625
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 19
$ws.Range("K26").Value = 26.666666666666
$ws.Range("L26").Value = 26.666666666666
$ws.Range("F27").Value = 5
$ws.Range("I27").Value = 31
$ws.Range("K27").Value = 3.333333333333
$ws.Range("L27").Value = 40.909090909090
